$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.93"
$ws.Range("E2").Value = "'1.93%"
$ws.Range("D3").Value = "'27.34"
$ws.Range("E3").Value = "'1.23%"
$ws.Range("D4").Value = "'4.719"
$ws.Range("E4").Value = "'8.27%"
$ws.Range("D5").Value = "'0.06072"
$ws.Range("E5").Value = "'3.16%"
$ws.Range("D6").Value = "'6.661"
$ws.Range("E6").Value = "'0.46%"
$ws.Range("D7").Value = "'0.8505"
$ws.Range("E7").Value = "'-0.14%"
$ws.Range("D8").Value = "'0.9166"
$ws.Range("E8").Value = "'-1.83%"
$ws.Range("D9").Value = "'0.1403"
$ws.Range("E9").Value = "'1.25%"
$ws.Range("D10").Value = "'0.04905"
$ws.Range("E10").Value = "'3.20%"
$ws.Range("D11").Value = "'0.07093"
$ws.Range("E11").Value = "'0.34%"
$ws.Range("D12").Value = "'0.03160"
$ws.Range("E12").Value = "'2.65%"
$ws.Range("D13").Value = "'0.09078"
$ws.Range("D14").Value = "'0.001530"
$ws.Range("E14").Value = "'0.18%"
$ws.Range("D15").Value = "'0.0006106"
$ws.Range("E15").Value = "'0.88%"
$ws.Range("D16").Value = "'0.006089"
$ws.Range("E16").Value = "'-0.92%"
$ws.Range("D17").Value = "'3.448"
$ws.Range("E17").Value = "'-1.03%"
$ws.Range("D18").Value = "'3.152"
$ws.Range("E18").Value = "'-0.52%"
$ws.Range("E19").Value = "'-1.27%"
$ws.Range("E20").Value = "'2.50%"
$ws.Range("D21").Value = "'0.1288"
$ws.Range("E21").Value = "'1.45%"
$ws.Range("D22").Value = "'4.092"
$ws.Range("E22").Value = "'4.49%"
$ws.Range("D23").Value = "'0.04242"
$ws.Range("E23").Value = "'-0.67%"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'0.09%"
$ws.Range("E25").Value = "'-8.55%"
$ws.Range("E26").Value = "'0.02%"
$ws.Range("E27").Value = "'3.07%"
$ws.Range("D40").Value = "'0.03871"
$ws.Range("E40").Value = "'1.63%"
$ws.Range("D41").Value = "'0.1114"
$ws.Range("E41").Value = "'1.38%"
$ws.Range("D42").Value = "'0.004127"
$ws.Range("D43").Value = "'0.01612"
$ws.Range("E43").Value = "'15.54%"
$ws.Range("D44").Value = "'0.002208"
$ws.Range("E44").Value = "'-9.51%"
$ws.Range("D45").Value = "'0.00005331"
$ws.Range("E45").Value = "'-0.93%"
$ws.Range("E46").Value = "'-0.03%"
$ws.Range("E47").Value = "'-4.30%"
$ws.Range("D48").Value = "'0.1321"
$ws.Range("E48").Value = "'-47.65%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.03%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.03%"
